# Update the build timestamp string wherever it appears in the workbook.
$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet: A2 (Version line) and A6 (Recommended Citation line) ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Shaqu No.1 Coal Mine, China, M0335, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# --- "Boundaries and methane sources" sheet: column S (build_version) rows 2-16 ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
